$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as literal text,
# even when the string looks like a number (e.g. "1.19", "0.692").
# Forcing the @ (text) number format before the write stops Excel
# from re-interpreting the string as a numeric value, and resetting
# the style back to Normal afterwards avoids leaving a stray
# direct-format on the cell.
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "44.224.74"
$ws.Range("E2").Value = "  +1.64%  "
Set-TextValue "D3" "2.365.05"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue "D5" "0.692"
$ws.Range("E5").Value = "  +5.36%  "
Set-TextValue "D6" "242.92"
$ws.Range("E6").Value = "  +2.89%  "
Set-TextValue "D7" "74.16"
$ws.Range("E7").Value = "  +2.93%  "
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue "D9" "0.603"
$ws.Range("E9").Value = "  +28.10%  "
$ws.Range("E10").Value = "  +6.49%  "
Set-TextValue "D11" "31.66"
$ws.Range("E11").Value = "  +16.70%  "
Set-TextValue "D12" "7.53"
$ws.Range("E12").Value = "  +19.86%  "
$ws.Range("E13").Value = "  +2.11%  "
Set-TextValue "D14" "2.716.80"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("E15").Value = "  +7.59%  "
$ws.Range("E16").Value = "  +7.27%  "
Set-TextValue "D17" "2.360.58"
$ws.Range("E17").Value = "  -0.82%  "
Set-TextValue "D18" "44.211.15"
$ws.Range("E18").Value = "  +1.67%  "
Set-TextValue "D19" "0.0000104"
$ws.Range("E19").Value = "  +4.69%  "
Set-TextValue "D20" "6.73"
$ws.Range("E20").Value = "  +5.52%  "
Set-TextValue "D21" "78.69"
$ws.Range("E21").Value = "  +5.28%  "
Set-TextValue "D22" "257.36"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("E25").Value = "  -3.09%  "
Set-TextValue "D26" "10.81"
$ws.Range("E26").Value = "  +8.31%  "
$ws.Range("E27").Value = "  +1.67%  "
Set-TextValue "D28" "1.66"
$ws.Range("E28").Value = "  +7.94%  "
Set-TextValue "D29" "22.74"
$ws.Range("E29").Value = "  -0.43%  "
Set-TextValue "D30" "175.22"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +6.69%  "
Set-TextValue "D33" "5.43"
$ws.Range("E33").Value = "  +8.51%  "
Set-TextValue "D34" "0.0761"
$ws.Range("E34").Value = "  +9.96%  "
Set-TextValue "D35" "5.41"
$ws.Range("E35").Value = "  +6.54%  "
Set-TextValue "D36" "3.92"
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("E39").Value = "  +7.65%  "
Set-TextValue "D40" "19.21"
$ws.Range("E40").Value = "  +0.70%  "
Set-TextValue "D41" "9.11"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").Value = "  +0.16%  "
Set-TextValue "D43" "0.200"
$ws.Range("E43").Value = "  +18.35%  "
$ws.Range("E44").Value = "  +5.25%  "
$ws.Range("E45").Value = "  +3.11%  "
Set-TextValue "D46" "2.50"
$ws.Range("E46").Value = "  +11.94%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D47" "1.19"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D48" "101.24"
$ws.Range("E48").Value = "  +1.47%  "
Set-TextValue "D49" "4.49"
$ws.Range("E49").Value = "  -1.11%  "
Set-TextValue "D50" "1.465.34"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue "D51" "0.000206"
$ws.Range("E51").Value = "  +2.09%  "
